$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values, entered in an order that reproduces the target shared-string table order ---
$ws.Range("A1").Value = "ECU CONTROLLER REQUIREMENTS"

$ws.Range("B3").Value = "Description"
$ws.Range("A3").Value = "Component "
$ws.Range("E3").Value = "Channel Requirement"
$ws.Range("A4").Value = "External Storage"
$ws.Range("D3").Value = "Quantity"
$ws.Range("A5").Value = "Solenoid Valve"
$ws.Range("C3").Value = "Type/Purpose"
$ws.Range("C4").Value = "Measurement"
$ws.Range("C5").Value = "Control"
$ws.Range("A6").Value = "Lithium ion Battery"
$ws.Range("C6").Value = "Utility"
$ws.Range("E5").Value = "GPIO"
$ws.Range("A7").Value = "Pump"
$ws.Range("C7").Value = "Control"
$ws.Range("A8").Value = "Spark igniter"
$ws.Range("C8").Value = "Control"

$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 2
$ws.Range("D6").Value = 1
$ws.Range("D7").Value = 2
$ws.Range("D8").Value = 1

# --- Title row formatting (row 1) ---
# Merge first so later per-sub-range formatting survives independently on each cell.
$ws.Range("A1:F1").Merge()

$titleFont = $ws.Range("A1:C1").Font
$titleFont.Bold = $true
$titleFont.Size = 16
$titleFont.Name = "Calibri (Body)"
$ws.Range("A1:C1").HorizontalAlignment = -4108

$ws.Range("D1:F1").HorizontalAlignment = -4108

$ws.Rows.Item(1).RowHeight = 21

# --- Column header row (row 3) formatting ---
$headerFont = $ws.Range("A3:E3").Font
$headerFont.Italic = $true
$headerFont.Underline = $true
$headerFont.Size = 14

$ws.Rows.Item(3).RowHeight = 19

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 16.166666666666668
$ws.Columns.Item(2).ColumnWidth = 13.76
$ws.Columns.Item(3).ColumnWidth = 13.76
$ws.Columns.Item(5).ColumnWidth = 18.166666666666668

# --- Selection matches post-edit cursor position ---
$ws.Range("A9").Select()

Write-Output "done"
